$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-Text "2025-07-15 Tuesday" "2025-07-16 Wednesday"
Replace-Text "498×4=" "906×6="
Replace-Text "840×8=" "984×6="
Replace-Text "247×3=" "381×3="
Replace-Text "664×5=" "337×7="
Replace-Text "333×4=" "146×5="
Replace-Text "821×2=" "682×3="
Replace-Text "598×8=" "229×3="
Replace-Text "984×3=" "342×6="
Replace-Text "444×5=" "200×8="
Replace-Text "410×8=" "884×9="
Replace-Text "547×2=" "860×7="
Replace-Text "377×5=" "272×9="
Replace-Text "207×7=" "984×3="
Replace-Text "193×7=" "978×4="
Replace-Text "844×7=" "977×5="
Replace-Text "308×3=" "832×5="
Replace-Text "418×2=" "517×5="
Replace-Text "120×2=" "936×5="
Replace-Text "172×5=" "489×7="
Replace-Text "686×7=" "917×2="
Replace-Text "126×2=" "265×4="
Replace-Text "953×2=" "470×5="
Replace-Text "642×5=" "386×5="
Replace-Text "651×3=" "270×7="
Replace-Text "350×7=" "169×2="
